$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 561.875
$ws.Cells.Item(2, 9).Value = 165.83333
$ws.Cells.Item(2, 10).Value = 1750
$ws.Cells.Item(2, 11).Value = 165.83333
$ws.Cells.Item(2, 12).Value = 1750
$ws.Cells.Item(2, 13).Value = -52.83332999999999
$ws.Cells.Item(2, 14).Value = -1976
$ws.Cells.Item(17, 8).Value = 969.9787
$ws.Cells.Item(17, 10).Value = 990
$ws.Cells.Item(17, 12).Value = 2970
$ws.Cells.Item(17, 14).Value = -3306
$ws.Cells.Item(69, 8).Value = 8624.277
$ws.Cells.Item(69, 9).Value = 5599.25
$ws.Cells.Item(69, 10).Value = 9488.571
$ws.Cells.Item(69, 11).Value = 16797.75
$ws.Cells.Item(69, 12).Value = 28465.713
$ws.Cells.Item(69, 13).Value = -15923.75
$ws.Cells.Item(69, 14).Value = -30213.713
$ws.Cells.Item(72, 8).Value = 8624.277
$ws.Cells.Item(72, 9).Value = 5599.25
$ws.Cells.Item(72, 10).Value = 9488.571
$ws.Cells.Item(72, 11).Value = 50393.25
$ws.Cells.Item(72, 12).Value = 85397.139
$ws.Cells.Item(72, 13).Value = -46025.25
$ws.Cells.Item(72, 14).Value = -94133.139
$ws.Cells.Item(113, 8).Value = 2759.75
$ws.Cells.Item(113, 10).Value = 2422.5
$ws.Cells.Item(113, 12).Value = 2422.5
$ws.Cells.Item(113, 14).Value = -8930.5
$ws.Cells.Item(132, 8).Value = 49677.55
$ws.Cells.Item(132, 9).Value = 54331.863
$ws.Cells.Item(132, 11).Value = 162995.589
$ws.Cells.Item(132, 13).Value = -160465.589
$ws.Cells.Item(137, 8).Value = 1353700.1
$ws.Cells.Item(137, 9).Value = 1147.6364
$ws.Cells.Item(137, 11).Value = 3442.9092
$ws.Cells.Item(137, 13).Value = -892.9092000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5851094
$ws.Cells.Item(32, 9).Value = 6292461.5
$ws.Cells.Item(32, 11).Value = 6292461.5
$ws.Cells.Item(32, 13).Value = -6292174.5
$ws.Cells.Item(45, 8).Value = 2792.8572
$ws.Cells.Item(45, 9).Value = 2891.6667
$ws.Cells.Item(45, 11).Value = 2891.6667
$ws.Cells.Item(45, 13).Value = -2514.6667
$ws.Cells.Item(61, 8).Value = 927617.25
$ws.Cells.Item(61, 9).Value = 981918.25
$ws.Cells.Item(61, 11).Value = 981918.25
$ws.Cells.Item(61, 13).Value = -981706.25
$ws.Cells.Item(110, 8).Value = 622.08
$ws.Cells.Item(110, 9).Value = 606.3333
$ws.Cells.Item(110, 11).Value = 606.3333
$ws.Cells.Item(110, 13).Value = 1438.6667
$ws.Cells.Item(122, 8).Value = 3191.125
$ws.Cells.Item(122, 9).Value = 3004.7827
$ws.Cells.Item(122, 11).Value = 9014.348100000001
$ws.Cells.Item(122, 13).Value = -6564.348100000001
$ws.Cells.Item(136, 8).Value = 927617.25
$ws.Cells.Item(136, 9).Value = 981918.25
$ws.Cells.Item(136, 11).Value = 2945754.75
$ws.Cells.Item(136, 13).Value = -2943204.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 992.6818
$ws.Cells.Item(94, 9).Value = 902.73334
$ws.Cells.Item(94, 10).Value = 1185.4286
$ws.Cells.Item(94, 11).Value = 902.73334
$ws.Cells.Item(94, 12).Value = 1185.4286
$ws.Cells.Item(94, 13).Value = -451.73334
$ws.Cells.Item(94, 14).Value = -2087.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 556
$ws.Cells.Item(10, 9).Value = 522.6
$ws.Cells.Item(10, 10).Value = 589.4
$ws.Cells.Item(10, 11).Value = 522.6
$ws.Cells.Item(10, 12).Value = 589.4
$ws.Cells.Item(10, 13).Value = -383.6
$ws.Cells.Item(10, 14).Value = -867.4
$ws.Cells.Item(21, 8).Value = 1200
$ws.Cells.Item(21, 9).Value = 800
$ws.Cells.Item(21, 11).Value = 800
$ws.Cells.Item(21, 13).Value = -565
$ws.Cells.Item(22, 8).Value = 859.7778
$ws.Cells.Item(22, 9).Value = 791.86664
$ws.Cells.Item(22, 10).Value = 1199.3334
$ws.Cells.Item(22, 11).Value = 791.86664
$ws.Cells.Item(22, 12).Value = 1199.3334
$ws.Cells.Item(22, 13).Value = -441.86664
$ws.Cells.Item(22, 14).Value = -1899.3334
$ws.Cells.Item(86, 8).Value = 6724.737
$ws.Cells.Item(86, 9).Value = 6078.385
$ws.Cells.Item(86, 10).Value = 8125.1665
$ws.Cells.Item(86, 11).Value = 6078.385
$ws.Cells.Item(86, 12).Value = 8125.1665
$ws.Cells.Item(86, 13).Value = -4955.385
$ws.Cells.Item(86, 14).Value = -10371.1665
$ws.Cells.Item(88, 8).Value = 38035
$ws.Cells.Item(88, 9).Value = 38035
$ws.Cells.Item(88, 11).Value = 38035
$ws.Cells.Item(88, 13).Value = -37629
$ws.Cells.Item(89, 8).Value = 6724.737
$ws.Cells.Item(89, 9).Value = 6078.385
$ws.Cells.Item(89, 10).Value = 8125.1665
$ws.Cells.Item(89, 11).Value = 30391.925
$ws.Cells.Item(89, 12).Value = 40625.8325
$ws.Cells.Item(89, 13).Value = -24775.925
$ws.Cells.Item(89, 14).Value = -51857.8325
$ws.Cells.Item(91, 8).Value = 38035
$ws.Cells.Item(91, 9).Value = 38035
$ws.Cells.Item(91, 11).Value = 38035
$ws.Cells.Item(91, 13).Value = -36631
$ws.Cells.Item(105, 8).Value = 26419.5
$ws.Cells.Item(105, 9).Value = 32852.184
$ws.Cells.Item(105, 10).Value = 2833
$ws.Cells.Item(105, 11).Value = 32852.184
$ws.Cells.Item(105, 12).Value = 2833
$ws.Cells.Item(105, 13).Value = -31105.184
$ws.Cells.Item(105, 14).Value = -6327
$ws.Cells.Item(134, 8).Value = 2278371
$ws.Cells.Item(134, 9).Value = 6490.269
$ws.Cells.Item(134, 10).Value = 10716786
$ws.Cells.Item(134, 11).Value = 19470.807
$ws.Cells.Item(134, 12).Value = 32150358
$ws.Cells.Item(134, 13).Value = -16935.807
$ws.Cells.Item(134, 14).Value = -32155428

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(93, 8).Value = 5814.1113
$ws.Cells.Item(93, 9).Value = 1999
$ws.Cells.Item(93, 10).Value = 6291
$ws.Cells.Item(93, 11).Value = 5997
$ws.Cells.Item(93, 12).Value = 18873
$ws.Cells.Item(93, 13).Value = -4125
$ws.Cells.Item(93, 14).Value = -22617
$ws.Cells.Item(107, 8).Value = 643.7143
$ws.Cells.Item(107, 9).Value = 500
$ws.Cells.Item(107, 10).Value = 701.2
$ws.Cells.Item(107, 11).Value = 1500
$ws.Cells.Item(107, 12).Value = 2103.6
$ws.Cells.Item(107, 13).Value = 420
$ws.Cells.Item(107, 14).Value = -5943.6
$ws.Cells.Item(132, 8).Value = 2198.6667
$ws.Cells.Item(132, 9).Value = 1289.1666
$ws.Cells.Item(132, 10).Value = 4017.6667
$ws.Cells.Item(132, 11).Value = 11602.4994
$ws.Cells.Item(132, 12).Value = 36159.0003
$ws.Cells.Item(132, 13).Value = -9072.499400000001
$ws.Cells.Item(132, 14).Value = -41219.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 23334.334
$ws.Cells.Item(40, 10).Value = 23334.334
$ws.Cells.Item(40, 12).Value = 23334.334
$ws.Cells.Item(40, 14).Value = -23636.334
$ws.Cells.Item(70, 8).Value = 7349.1665
$ws.Cells.Item(70, 9).Value = 7349.1665
$ws.Cells.Item(70, 11).Value = 7349.1665
$ws.Cells.Item(70, 13).Value = -7079.1665
$ws.Cells.Item(73, 8).Value = 7349.1665
$ws.Cells.Item(73, 9).Value = 7349.1665
$ws.Cells.Item(73, 11).Value = 7349.1665
$ws.Cells.Item(73, 13).Value = -6413.1665
$ws.Cells.Item(102, 8).Value = 3542.225
$ws.Cells.Item(102, 9).Value = 2976.8857
$ws.Cells.Item(102, 10).Value = 7499.6
$ws.Cells.Item(102, 11).Value = 2976.8857
$ws.Cells.Item(102, 12).Value = 7499.6
$ws.Cells.Item(102, 13).Value = -1354.8857
$ws.Cells.Item(102, 14).Value = -10743.6
$ws.Cells.Item(107, 8).Value = 28365.385
$ws.Cells.Item(107, 9).Value = 44643.5
$ws.Cells.Item(107, 11).Value = 44643.5
$ws.Cells.Item(107, 13).Value = -42723.5
$ws.Cells.Item(122, 8).Value = 27651.418
$ws.Cells.Item(122, 9).Value = 48046.316
$ws.Cells.Item(122, 10).Value = 6285.3335
$ws.Cells.Item(122, 11).Value = 144138.948
$ws.Cells.Item(122, 12).Value = 18856.0005
$ws.Cells.Item(122, 13).Value = -141688.948
$ws.Cells.Item(122, 14).Value = -23756.0005
$ws.Cells.Item(132, 8).Value = 28923630
$ws.Cells.Item(132, 9).Value = 37490948
$ws.Cells.Item(132, 11).Value = 112472844
$ws.Cells.Item(132, 13).Value = -112470314

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 7700.467
$ws.Cells.Item(16, 9).Value = 8716
$ws.Cells.Item(16, 11).Value = 8716
$ws.Cells.Item(16, 13).Value = -8546
$ws.Cells.Item(22, 8).Value = 522.7
$ws.Cells.Item(22, 9).Value = 329.5
$ws.Cells.Item(22, 10).Value = 812.5
$ws.Cells.Item(22, 11).Value = 329.5
$ws.Cells.Item(22, 12).Value = 812.5
$ws.Cells.Item(22, 13).Value = -34.5
$ws.Cells.Item(22, 14).Value = -1402.5
$ws.Cells.Item(27, 8).Value = 522.7
$ws.Cells.Item(27, 9).Value = 329.5
$ws.Cells.Item(27, 10).Value = 812.5
$ws.Cells.Item(27, 11).Value = 329.5
$ws.Cells.Item(27, 12).Value = 812.5
$ws.Cells.Item(27, 13).Value = -222.5
$ws.Cells.Item(27, 14).Value = -1026.5
$ws.Cells.Item(40, 8).Value = 3507.7222
$ws.Cells.Item(40, 9).Value = 3250.4546
$ws.Cells.Item(40, 11).Value = 3250.4546
$ws.Cells.Item(40, 13).Value = -3114.4546
$ws.Cells.Item(43, 8).Value = 295999.84
$ws.Cells.Item(43, 9).Value = 6000
$ws.Cells.Item(43, 11).Value = 6000
$ws.Cells.Item(43, 13).Value = -5807
$ws.Cells.Item(55, 8).Value = 15625231
$ws.Cells.Item(55, 9).Value = 259.77777
$ws.Cells.Item(55, 11).Value = 259.77777
$ws.Cells.Item(55, 13).Value = -86.77776999999998
$ws.Cells.Item(82, 8).Value = 1079.2858
$ws.Cells.Item(82, 10).Value = 1156.5
$ws.Cells.Item(82, 12).Value = 1156.5
$ws.Cells.Item(82, 14).Value = -1878.5
$ws.Cells.Item(85, 8).Value = 1079.2858
$ws.Cells.Item(85, 10).Value = 1156.5
$ws.Cells.Item(85, 12).Value = 1156.5
$ws.Cells.Item(85, 14).Value = -3652.5
$ws.Cells.Item(122, 8).Value = 5520
$ws.Cells.Item(122, 9).Value = 5233.3335
$ws.Cells.Item(122, 10).Value = 6666.6665
$ws.Cells.Item(122, 11).Value = 15700.0005
$ws.Cells.Item(122, 12).Value = 19999.9995
$ws.Cells.Item(122, 13).Value = -13250.0005
$ws.Cells.Item(122, 14).Value = -24899.9995
$ws.Cells.Item(132, 8).Value = 647772.3
$ws.Cells.Item(132, 9).Value = 968067.1
$ws.Cells.Item(132, 10).Value = 7182.6113
$ws.Cells.Item(132, 11).Value = 2904201.3
$ws.Cells.Item(132, 12).Value = 21547.8339
$ws.Cells.Item(132, 13).Value = -2901671.3
$ws.Cells.Item(132, 14).Value = -26607.8339

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2364.8445
$ws.Cells.Item(122, 9).Value = 2026.6471
$ws.Cells.Item(122, 10).Value = 3410.182
$ws.Cells.Item(122, 11).Value = 6079.9413
$ws.Cells.Item(122, 12).Value = 10230.546
$ws.Cells.Item(122, 13).Value = -3629.9413
$ws.Cells.Item(122, 14).Value = -15130.546
$ws.Cells.Item(132, 8).Value = 10602
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
